$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2952119588580189
$ws.Range("D2").Value = 0.03072992620381854
$ws.Range("E2").Value = 0.0830438487301226
$ws.Range("F2").Value = 8.821586488591208
$ws.Range("G2").Value = 0.002784564261579222
$ws.Range("J2").Value = 0.3517665580556226
$ws.Range("K2").Value = 4.900726431097041
$ws.Range("L2").Value = 0.05085815343876732
$ws.Range("M2").Value = 1.135977361395618
$ws.Range("C3").Value = 0.2948735464759125
$ws.Range("D3").Value = 0.02766620663722819
$ws.Range("E3").Value = 0.08357188095053836
$ws.Range("F3").Value = 8.677735975671226
$ws.Range("G3").Value = 0.0027919279193664
$ws.Range("J3").Value = 0.3499071365648376
$ws.Range("K3").Value = 4.858164982595042
$ws.Range("L3").Value = 0.0507789119186004
$ws.Range("M3").Value = 1.132013603506017
$ws.Range("C4").Value = 0.2948059669456597
$ws.Range("D4").Value = 0.0257755073485626
$ws.Range("E4").Value = 0.08392254424038548
$ws.Range("F4").Value = 8.59162570014243
$ws.Range("G4").Value = 0.002796681576330024
$ws.Range("J4").Value = 0.3488366179803464
$ws.Range("K4").Value = 4.835765742258872
$ws.Range("L4").Value = 0.05073243591338894
$ws.Range("M4").Value = 1.130365205640459
$ws.Range("C5").Value = 0.2948136660047709
$ws.Range("D5").Value = 0.0250024329329861
$ws.Range("E5").Value = 0.08407211073835352
$ws.Range("F5").Value = 8.557084902218662
$ws.Range("G5").Value = 0.002798677380747511
$ws.Range("J5").Value = 0.3484181559830475
$ws.Range("K5").Value = 4.827574935387759
$ws.Range("L5").Value = 0.05071403927312801
$ws.Range("M5").Value = 1.129890765341983
$ws.Range("C6").Value = 0.2948170722604431
$ws.Range("D6").Value = 0.02487390186446703
$ws.Range("E6").Value = 0.08409734940949853
$ws.Range("F6").Value = 8.551382434076288
$ws.Range("G6").Value = 0.002799012331344322
$ws.Range("J6").Value = 0.3483497414304395
$ws.Range("K6").Value = 4.826271404022975
$ws.Range("L6").Value = 0.05071101715505932
$ws.Range("M6").Value = 1.129823895788199
$ws.Range("C7").Value = 0.2948059281203257
$ws.Range("D7").Value = 0.02576509217820444
$ws.Range("E7").Value = 0.08392453432292779
$ws.Range("F7").Value = 8.591157653801247
$ws.Range("G7").Value = 0.002796708254615917
$ws.Range("J7").Value = 0.3488309026017262
$ws.Range("K7").Value = 4.835651486300975
$ws.Range("L7").Value = 0.05073218561916804
$ws.Range("M7").Value = 1.130358008565167
$ws.Range("C8").Value = 0.2950661619270534
$ws.Range("D8").Value = 0.02967542264926948
$ws.Range("E8").Value = 0.0832204350983563
$ws.Range("F8").Value = 8.771522912636641
$ws.Range("G8").Value = 0.002787055165394818
$ws.Range("J8").Value = 0.3511105874350378
$ws.Range("K8").Value = 4.885275182906128
$ws.Range("L8").Value = 0.05083037565997195
$ws.Range("M8").Value = 1.134447529464438
$ws.Range("C9").Value = 0.2966902818171349
$ws.Range("D9").Value = 0.03727672916865288
$ws.Range("E9").Value = 0.08204880032660178
$ws.Range("F9").Value = 9.143129458360079
$ws.Range("G9").Value = 0.002769958641556208
$ws.Range("J9").Value = 0.3561513393431355
$ws.Range("K9").Value = 5.01231154001988
$ws.Range("L9").Value = 0.05104045903180499
$ws.Range("M9").Value = 1.148710487552371
$ws.Range("C10").Value = 0.2985651070582094
$ws.Range("D10").Value = 0.04283357139823352
$ws.Range("E10").Value = 0.08131445611747079
$ws.Range("F10").Value = 9.427571210749932
$ws.Range("G10").Value = 0.002758500813565583
$ws.Range("J10").Value = 0.36021075785996
$ws.Range("K10").Value = 5.123930884425477
$ws.Range("L10").Value = 0.05120586051727294
$ws.Range("M10").Value = 1.163016884807448
$ws.Range("C11").Value = 0.299566666619441
$ws.Range("D11").Value = 0.04535830069586666
$ws.Range("E11").Value = 0.08100763111838027
$ws.Range("F11").Value = 9.559563492971108
$ws.Range("G11").Value = 0.00275352473223927
$ws.Range("J11").Value = 0.3621366494210108
$ws.Range("K11").Value = 5.17871911236432
$ws.Range("L11").Value = 0.05128358751416862
$ws.Range("M11").Value = 1.170361544873053
$ws.Range("C12").Value = 0.2999673602179627
$ws.Range("D12").Value = 0.04631412608091523
$ws.Range("E12").Value = 0.08089534376408913
$ws.Range("F12").Value = 9.609927815768231
$ws.Range("G12").Value = 0.00275167413505485
$ws.Range("J12").Value = 0.362877473159557
$ws.Range("K12").Value = 5.200045876312572
$ws.Range("L12").Value = 0.05131338443410449
$ws.Range("M12").Value = 1.173263471245875
$ws.Range("C13").Value = 0.2998801100928148
$ws.Range("D13").Value = 0.04610827903232462
$ws.Range("E13").Value = 0.08091935361516711
$ws.Range("F13").Value = 9.599063850996004
$ws.Range("G13").Value = 0.002752071197300986
$ws.Range("J13").Value = 0.3627174085103348
$ws.Range("K13").Value = 5.195426951164166
$ws.Range("L13").Value = 0.05130695087458115
$ws.Range("M13").Value = 1.172633116981942
$ws.Range("C14").Value = 0.2995992022824367
$ws.Range("D14").Value = 0.04543694051019997
$ws.Range("E14").Value = 0.08099831508287636
$ws.Range("F14").Value = 9.56369929934408
$ws.Range("G14").Value = 0.002753371807486196
$ws.Range("J14").Value = 0.3621973655637731
$ws.Range("K14").Value = 5.180462044926969
$ws.Range("L14").Value = 0.05128603160671119
$ws.Range("M14").Value = 1.170597867694319
$ws.Range("C15").Value = 0.2994299296002509
$ws.Range("D15").Value = 0.04502570202573963
$ws.Range("E15").Value = 0.08104718878709427
$ws.Range("F15").Value = 9.542087442233992
$ws.Range("G15").Value = 0.002754172855659547
$ws.Range("J15").Value = 0.361880329908864
$ws.Range("K15").Value = 5.171371178541904
$ws.Range("L15").Value = 0.05127326545674116
$ws.Range("M15").Value = 1.169366943691088
$ws.Range("C16").Value = 0.2985026484504658
$ws.Range("D16").Value = 0.04266852996258308
$ws.Range("E16").Value = 0.08133505451447753
$ws.Range("F16").Value = 9.418998176650291
$ws.Range("G16").Value = 0.002758830745297924
$ws.Range("J16").Value = 0.3600865023924626
$ws.Range("K16").Value = 5.120431307265733
$ws.Range("L16").Value = 0.05120083139466658
$ws.Range("M16").Value = 1.162553758926819
$ws.Range("C17").Value = 0.2979719046875289
$ws.Range("D17").Value = 0.04122183581051786
$ws.Range("E17").Value = 0.08151861458901788
$ws.Range("F17").Value = 9.344158092391524
$ws.Range("G17").Value = 0.002761748534926653
$ws.Range("J17").Value = 0.3590064386630587
$ws.Range("K17").Value = 5.090210822083463
$ws.Range("L17").Value = 0.05115703618679035
$ws.Range("M17").Value = 1.158588614932746
$ws.Range("C18").Value = 0.297680627755966
$ws.Range("D18").Value = 0.04038943045014776
$ws.Range("E18").Value = 0.08162675758117288
$ws.Range("F18").Value = 9.301356129051243
$ws.Range("G18").Value = 0.00276344901023684
$ws.Range("J18").Value = 0.3583926649225262
$ws.Range("K18").Value = 5.073206258599384
$ws.Range("L18").Value = 0.05113207987550616
$ws.Range("M18").Value = 1.156386691608546
$ws.Range("C19").Value = 0.2975844084185724
$ws.Range("D19").Value = 0.04010753346847196
$ws.Range("E19").Value = 0.08166381379504983
$ws.Range("F19").Value = 9.28690579557383
$ws.Range("G19").Value = 0.00276402858847654
$ws.Range("J19").Value = 0.3581861271967881
$ws.Range("K19").Value = 5.067513553718868
$ws.Range("L19").Value = 0.05112367002480767
$ws.Range("M19").Value = 1.155654668290651
$ws.Range("C20").Value = 0.2980269547597203
$ws.Range("D20").Value = 0.04137586892625222
$ws.Range("E20").Value = 0.08149880904363815
$ws.Range("F20").Value = 9.352099628339829
$ws.Range("G20").Value = 0.002761435630893259
$ws.Range("J20").Value = 0.3591206412342487
$ws.Range("K20").Value = 5.093388759286654
$ws.Range("L20").Value = 0.05116167404782423
$ws.Range("M20").Value = 1.159002561256131
$ws.Range("C21").Value = 0.2996811298965696
$ws.Range("D21").Value = 0.04563413358231116
$ws.Range("E21").Value = 0.08097501644482108
$ws.Range("F21").Value = 9.574076300872491
$ws.Range("G21").Value = 0.002752988872758264
$ws.Range("J21").Value = 0.3623498007590058
$ws.Range("K21").Value = 5.184841846707002
$ws.Range("L21").Value = 0.05129216619146959
$ws.Range("M21").Value = 1.171192392003931
$ws.Range("C22").Value = 0.3008871309513097
$ws.Range("D22").Value = 0.04841588686231546
$ws.Range("E22").Value = 0.08065541653647834
$ws.Range("F22").Value = 9.721379631309503
$ws.Range("G22").Value = 0.002747664972785179
$ws.Range("J22").Value = 0.3645275159265253
$ws.Range("K22").Value = 5.247991733870833
$ws.Range("L22").Value = 0.05137957137679194
$ws.Range("M22").Value = 1.179862650798341
$ws.Range("C23").Value = 0.3002320219034686
$ws.Range("D23").Value = 0.04693125835228784
$ws.Range("E23").Value = 0.08082391832371094
$ws.Range("F23").Value = 9.642554437785122
$ws.Range("G23").Value = 0.002750488527088848
$ws.Range("J23").Value = 0.3633590274451919
$ws.Range("K23").Value = 5.213977245096316
$ws.Range("L23").Value = 0.05133272548128076
$ws.Range("M23").Value = 1.175170675355616
$ws.Range("C24").Value = 0.2980020234697776
$ws.Range("D24").Value = 0.04130623269404055
$ws.Range("E24").Value = 0.081507754994254
$ws.Range("F24").Value = 9.348508560974039
$ws.Range("G24").Value = 0.00276157702318004
$ws.Range("J24").Value = 0.3590689879314937
$ws.Range("K24").Value = 5.091950862872409
$ws.Range("L24").Value = 0.05115957657927694
$ws.Range("M24").Value = 1.158815174160772
$ws.Range("C25").Value = 0.2961314183961292
$ws.Range("D25").Value = 0.0352262015897935
$ws.Range("E25").Value = 0.08234347891389326
$ws.Range("F25").Value = 9.040630941947796
$ws.Range("G25").Value = 0.002774388957803986
$ws.Range("J25").Value = 0.3547257773397376
$ws.Range("K25").Value = 4.974746022191255
$ws.Range("L25").Value = 0.0183191443987436
$ws.Range("M25").Value = 1.144181460155359
